$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 3918.3333
$ws.Cells.Item(15, 9).Value = 3918.3333
$ws.Cells.Item(15, 11).Value = 11754.9999
$ws.Cells.Item(15, 13).Value = -11585.9999
$ws.Cells.Item(98, 8).Value = 5281.3335
$ws.Cells.Item(98, 9).Value = 4337.6
$ws.Cells.Item(98, 10).Value = 10000
$ws.Cells.Item(98, 11).Value = 4337.6
$ws.Cells.Item(98, 12).Value = 10000
$ws.Cells.Item(98, 13).Value = -2839.6
$ws.Cells.Item(98, 14).Value = -12996
$ws.Cells.Item(112, 8).Value = 5001234
$ws.Cells.Item(112, 9).Value = 2057.1428
$ws.Cells.Item(112, 10).Value = 5815053.5
$ws.Cells.Item(112, 11).Value = 6171.428400000001
$ws.Cells.Item(112, 12).Value = 17445160.5
$ws.Cells.Item(112, 13).Value = -5063.428400000001
$ws.Cells.Item(112, 14).Value = -17447376.5
$ws.Cells.Item(122, 8).Value = 5281.3335
$ws.Cells.Item(122, 9).Value = 4337.6
$ws.Cells.Item(122, 10).Value = 10000
$ws.Cells.Item(122, 11).Value = 13012.8
$ws.Cells.Item(122, 12).Value = 30000
$ws.Cells.Item(122, 13).Value = -10562.8
$ws.Cells.Item(122, 14).Value = -34900
$ws.Cells.Item(129, 8).Value = 3248101.8
$ws.Cells.Item(129, 10).Value = 1492.9207
$ws.Cells.Item(129, 12).Value = 4478.7621
$ws.Cells.Item(129, 14).Value = -14478.7621
$ws.Cells.Item(141, 8).Value = 1070459.8
$ws.Cells.Item(141, 9).Value = 1969.2632
$ws.Cells.Item(141, 10).Value = 5130724
$ws.Cells.Item(141, 11).Value = 5907.7896
$ws.Cells.Item(141, 12).Value = 15392172
$ws.Cells.Item(141, 13).Value = -727.7896000000001
$ws.Cells.Item(141, 14).Value = -15402532

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2006.32
$ws.Cells.Item(61, 9).Value = 1462.45
$ws.Cells.Item(61, 10).Value = 4181.8
$ws.Cells.Item(61, 11).Value = 1462.45
$ws.Cells.Item(61, 12).Value = 4181.8
$ws.Cells.Item(61, 13).Value = -1250.45
$ws.Cells.Item(61, 14).Value = -4605.8
$ws.Cells.Item(74, 8).Value = 869.6923
$ws.Cells.Item(74, 9).Value = 775.5
$ws.Cells.Item(74, 10).Value = 2000
$ws.Cells.Item(74, 11).Value = 775.5
$ws.Cells.Item(74, 12).Value = 2000
$ws.Cells.Item(74, 13).Value = 98.5
$ws.Cells.Item(74, 14).Value = -3748
$ws.Cells.Item(77, 8).Value = 869.6923
$ws.Cells.Item(77, 9).Value = 775.5
$ws.Cells.Item(77, 10).Value = 2000
$ws.Cells.Item(77, 11).Value = 3877.5
$ws.Cells.Item(77, 12).Value = 10000
$ws.Cells.Item(77, 13).Value = 490.5
$ws.Cells.Item(77, 14).Value = -18736
$ws.Cells.Item(132, 8).Value = 27782754
$ws.Cells.Item(132, 9).Value = 41671616
$ws.Cells.Item(132, 11).Value = 125014848
$ws.Cells.Item(132, 13).Value = -125012318
$ws.Cells.Item(136, 8).Value = 2006.32
$ws.Cells.Item(136, 9).Value = 1462.45
$ws.Cells.Item(136, 10).Value = 4181.8
$ws.Cells.Item(136, 11).Value = 4387.35
$ws.Cells.Item(136, 12).Value = 12545.4
$ws.Cells.Item(136, 13).Value = -1837.35
$ws.Cells.Item(136, 14).Value = -17645.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3779.5833
$ws.Cells.Item(134, 9).Value = 3214.0908
$ws.Cells.Item(134, 10).Value = 10000
$ws.Cells.Item(134, 11).Value = 9642.2724
$ws.Cells.Item(134, 12).Value = 30000
$ws.Cells.Item(134, 13).Value = -7107.2724
$ws.Cells.Item(134, 14).Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3033361.5
$ws.Cells.Item(31, 9).Value = 4002417.2
$ws.Cells.Item(31, 10).Value = 5062.5
$ws.Cells.Item(31, 11).Value = 4002417.2
$ws.Cells.Item(31, 12).Value = 5062.5
$ws.Cells.Item(31, 13).Value = -4002122.2
$ws.Cells.Item(31, 14).Value = -5652.5
$ws.Cells.Item(34, 8).Value = 3033361.5
$ws.Cells.Item(34, 9).Value = 4002417.2
$ws.Cells.Item(34, 10).Value = 5062.5
$ws.Cells.Item(34, 11).Value = 4002417.2
$ws.Cells.Item(34, 12).Value = 5062.5
$ws.Cells.Item(34, 13).Value = -4002215.2
$ws.Cells.Item(34, 14).Value = -5466.5
$ws.Cells.Item(58, 8).Value = 27781448
$ws.Cells.Item(58, 9).Value = 1246.6666
$ws.Cells.Item(58, 10).Value = 41671550
$ws.Cells.Item(58, 11).Value = 1246.6666
$ws.Cells.Item(58, 12).Value = 41671550
$ws.Cells.Item(58, 13).Value = -1043.6666
$ws.Cells.Item(58, 14).Value = -41671956
$ws.Cells.Item(132, 8).Value = 3909.25
$ws.Cells.Item(132, 9).Value = 2748.8333
$ws.Cells.Item(132, 11).Value = 8246.499899999999
$ws.Cells.Item(132, 13).Value = -5716.499899999999
$ws.Cells.Item(134, 8).Value = 1788.4
$ws.Cells.Item(134, 9).Value = 939.46155
$ws.Cells.Item(134, 10).Value = 3365
$ws.Cells.Item(134, 11).Value = 2818.38465
$ws.Cells.Item(134, 12).Value = 10095
$ws.Cells.Item(134, 13).Value = -283.38465
$ws.Cells.Item(134, 14).Value = -15165
$ws.Cells.Item(136, 8).Value = 27781448
$ws.Cells.Item(136, 9).Value = 1246.6666
$ws.Cells.Item(136, 10).Value = 41671550
$ws.Cells.Item(136, 11).Value = 3739.9998
$ws.Cells.Item(136, 12).Value = 125014650
$ws.Cells.Item(136, 13).Value = -1189.9998
$ws.Cells.Item(136, 14).Value = -125019750

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 134
$ws.Cells.Item(23, 9).Value = 78
$ws.Cells.Item(23, 10).Value = 142
$ws.Cells.Item(23, 11).Value = 234
$ws.Cells.Item(23, 12).Value = 426
$ws.Cells.Item(23, 13).Value = 1
$ws.Cells.Item(23, 14).Value = -896
$ws.Cells.Item(58, 8).Value = 1716.6666
$ws.Cells.Item(58, 9).Value = 500
$ws.Cells.Item(58, 10).Value = 2933.3333
$ws.Cells.Item(58, 11).Value = 1500
$ws.Cells.Item(58, 12).Value = 8799.999899999999
$ws.Cells.Item(58, 13).Value = -1372
$ws.Cells.Item(58, 14).Value = -9055.999899999999
$ws.Cells.Item(132, 8).Value = 3570.2856
$ws.Cells.Item(132, 9).Value = 798
$ws.Cells.Item(132, 10).Value = 7266.6665
$ws.Cells.Item(132, 11).Value = 7182
$ws.Cells.Item(132, 12).Value = 65399.9985
$ws.Cells.Item(132, 13).Value = -4652
$ws.Cells.Item(132, 14).Value = -70459.9985
$ws.Cells.Item(133, 8).Value = 6756.5557
$ws.Cells.Item(133, 9).Value = 12343.333
$ws.Cells.Item(133, 10).Value = 3963.1667
$ws.Cells.Item(133, 11).Value = 37029.999
$ws.Cells.Item(133, 12).Value = 11889.5001
$ws.Cells.Item(133, 13).Value = -31969.999
$ws.Cells.Item(133, 14).Value = -22009.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3406.5483
$ws.Cells.Item(132, 9).Value = 2882.647
$ws.Cells.Item(132, 10).Value = 4042.7144
$ws.Cells.Item(132, 11).Value = 8647.940999999999
$ws.Cells.Item(132, 12).Value = 12128.1432
$ws.Cells.Item(132, 13).Value = -6117.940999999999
$ws.Cells.Item(132, 14).Value = -17188.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2898.4285
$ws.Cells.Item(82, 9).Value = 2175.3333
$ws.Cells.Item(82, 11).Value = 2175.3333
$ws.Cells.Item(82, 13).Value = -1814.3333
$ws.Cells.Item(85, 8).Value = 2898.4285
$ws.Cells.Item(85, 9).Value = 2175.3333
$ws.Cells.Item(85, 11).Value = 2175.3333
$ws.Cells.Item(85, 13).Value = -927.3332999999998
$ws.Cells.Item(132, 8).Value = 2781.4324
$ws.Cells.Item(132, 9).Value = 1548.0952
$ws.Cells.Item(132, 10).Value = 4400.1875
$ws.Cells.Item(132, 11).Value = 4644.2856
$ws.Cells.Item(132, 12).Value = 13200.5625
$ws.Cells.Item(132, 13).Value = -2114.2856
$ws.Cells.Item(132, 14).Value = -18260.5625
$ws.Cells.Item(136, 8).Value = 3127141
$ws.Cells.Item(136, 9).Value = 5264733
$ws.Cells.Item(136, 10).Value = 2968.2307
$ws.Cells.Item(136, 11).Value = 15794199
$ws.Cells.Item(136, 12).Value = 8904.6921
$ws.Cells.Item(136, 13).Value = -15791649
$ws.Cells.Item(136, 14).Value = -14004.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 2500
$ws.Cells.Item(62, 10).Value = 3000
$ws.Cells.Item(62, 12).Value = 3000
$ws.Cells.Item(62, 14).Value = -4248
$ws.Cells.Item(65, 8).Value = 2500
$ws.Cells.Item(65, 10).Value = 3000
$ws.Cells.Item(65, 12).Value = 15000
$ws.Cells.Item(65, 14).Value = -21240
$ws.Cells.Item(132, 8).Value = 189554.2
$ws.Cells.Item(132, 9).Value = 252140.97
$ws.Cells.Item(132, 10).Value = 10734.857
$ws.Cells.Item(132, 11).Value = 756422.91
$ws.Cells.Item(132, 12).Value = 32204.571
$ws.Cells.Item(132, 13).Value = -753892.91
$ws.Cells.Item(132, 14).Value = -37264.571
$ws.Cells.Item(136, 8).Value = 2044.1428
$ws.Cells.Item(136, 9).Value = 979.75
$ws.Cells.Item(136, 10).Value = 3463.3333
$ws.Cells.Item(136, 11).Value = 2939.25
$ws.Cells.Item(136, 12).Value = 10389.9999
$ws.Cells.Item(136, 13).Value = -389.25
$ws.Cells.Item(136, 14).Value = -15489.9999
